$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 187.5
$ws.Range("J12").Value = 225
$ws.Range("L12").Value = 225
$ws.Range("N12").Value = -565
$ws.Range("H18").Value = 1131.25
$ws.Range("I18").Value = 1042.8572
$ws.Range("K18").Value = 1042.8572
$ws.Range("M18").Value = -758.8571999999999
$ws.Range("H51").Value = 116759.6
$ws.Range("I51").Value = 9999
$ws.Range("J51").Value = 143449.75
$ws.Range("K51").Value = 9999
$ws.Range("L51").Value = 143449.75
$ws.Range("M51").Value = -9515
$ws.Range("N51").Value = -144417.75
$ws.Range("H74").Value = 2903
$ws.Range("I74").Value = 2903
$ws.Range("K74").Value = 2903
$ws.Range("M74").Value = -1967
$ws.Range("H77").Value = 2903
$ws.Range("I77").Value = 2903
$ws.Range("K77").Value = 14515
$ws.Range("M77").Value = -9835
$ws.Range("H80").Value = 1600
$ws.Range("H83").Value = 1600
$ws.Range("H86").Value = 11857.714
$ws.Range("I86").Value = 5000
$ws.Range("J86").Value = 13000.667
$ws.Range("K86").Value = 5000
$ws.Range("L86").Value = 13000.667
$ws.Range("M86").Value = -3877
$ws.Range("N86").Value = -15246.667
$ws.Range("H89").Value = 11857.714
$ws.Range("I89").Value = 5000
$ws.Range("J89").Value = 13000.667
$ws.Range("K89").Value = 25000
$ws.Range("L89").Value = 65003.335
$ws.Range("M89").Value = -19384
$ws.Range("N89").Value = -76235.33499999999
$ws.Range("H92").Value = 301
$ws.Range("I92").Value = 318.75
$ws.Range("J92").Value = 230
$ws.Range("K92").Value = 318.75
$ws.Range("L92").Value = 230
$ws.Range("M92").Value = 929.25
$ws.Range("N92").Value = -2726
$ws.Range("H137").Value = 2831.3333
$ws.Range("J137").Value = 2999.5
$ws.Range("L137").Value = 8998.5
$ws.Range("N137").Value = -14098.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1179.8
$ws.Range("I32").Value = 974.75
$ws.Range("K32").Value = 974.75
$ws.Range("M32").Value = -687.75
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").Value = $null
$ws.Range("H50").Value = 22666.334
$ws.Range("J50").Value = 23190
$ws.Range("L50").Value = 23190
$ws.Range("N50").Value = -24618
$ws.Range("H88").Value = 3207.889
$ws.Range("I88").Value = 699.5
$ws.Range("J88").Value = 3924.5715
$ws.Range("K88").Value = 699.5
$ws.Range("L88").Value = 3924.5715
$ws.Range("M88").Value = -293.5
$ws.Range("N88").Value = -4736.5715
$ws.Range("H91").Value = 3207.889
$ws.Range("I91").Value = 699.5
$ws.Range("J91").Value = 3924.5715
$ws.Range("K91").Value = 699.5
$ws.Range("L91").Value = 3924.5715
$ws.Range("M91").Value = 704.5
$ws.Range("N91").Value = -6732.5715
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").Value = $null
$ws.Range("H134").Value = 3747.5
$ws.Range("I134").Value = 3747.5
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 11242.5
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -8707.5
$ws.Range("N134").Value = $null
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 770.5833
$ws.Range("I31").Value = 770.5833
$ws.Range("K31").Value = 770.5833
$ws.Range("M31").Value = -475.5833
$ws.Range("H34").Value = 770.5833
$ws.Range("I34").Value = 770.5833
$ws.Range("K34").Value = 770.5833
$ws.Range("M34").Value = -568.5833
$ws.Range("H69").Value = 8000
$ws.Range("I69").Value = 8000
$ws.Range("K69").Value = 8000
$ws.Range("M69").Value = -7251
$ws.Range("H72").Value = 8000
$ws.Range("I72").Value = 8000
$ws.Range("K72").Value = 24000
$ws.Range("M72").Value = -20256
$ws.Range("H132").Value = 1478.75
$ws.Range("I132").Value = 1478.75
$ws.Range("K132").Value = 4436.25
$ws.Range("M132").Value = -1906.25
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 1006.75
$ws.Range("J41").Value = 1009.3333
$ws.Range("L41").Value = 3027.9999
$ws.Range("N41").Value = -3703.9999
$ws.Range("H50").Value = 602.2
$ws.Range("I50").Value = 168.33333
$ws.Range("J50").Value = 1253
$ws.Range("K50").Value = 504.99999
$ws.Range("L50").Value = 3759
$ws.Range("M50").Value = -23.99998999999997
$ws.Range("N50").Value = -4721
$ws.Range("H53").Value = 602.2
$ws.Range("I53").Value = 168.33333
$ws.Range("J53").Value = 1253
$ws.Range("K53").Value = 504.99999
$ws.Range("L53").Value = 3759
$ws.Range("M53").Value = -23.99998999999997
$ws.Range("N53").Value = -4721
$ws.Range("H69").Value = 700
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").Value = $null
$ws.Range("H72").Value = 700
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").Value = $null
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3000
$ws.Range("I132").Value = 3000
$ws.Range("K132").Value = 9000
$ws.Range("M132").Value = -6470
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4831.3335
$ws.Range("I61").Value = 4831.3335
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 4831.3335
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -4629.3335
$ws.Range("N61").Value = $null
$ws.Range("H93").Value = 2343.7778
$ws.Range("I93").Value = 2156.4285
$ws.Range("J93").Value = 2999.5
$ws.Range("K93").Value = 2156.4285
$ws.Range("L93").Value = 2999.5
$ws.Range("M93").Value = -908.4285
$ws.Range("N93").Value = -5495.5
$ws.Range("H113").Value = 4831.3335
$ws.Range("I113").Value = 4831.3335
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 4831.3335
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -2661.3335
$ws.Range("N113").Value = $null
$ws.Range("H133").Value = 120000
$ws.Range("J133").Value = 120000
$ws.Range("L133").Value = 120000
$ws.Range("N133").Value = -125060
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").Value = $null
$ws.Range("H105").Value = 30000
$ws.Range("J105").Value = 30000
$ws.Range("L105").Value = 30000
$ws.Range("N105").Value = -36988
$ws.Range("H132").Value = 1079.8334
$ws.Range("I132").Value = 1079.8334
$ws.Range("K132").Value = 3239.5002
$ws.Range("M132").Value = -709.5001999999999
$ws.Range("H136").Value = 19998.334
$ws.Range("I136").Value = 19998.334
$ws.Range("K136").Value = 59995.00199999999
$ws.Range("M136").Value = -57445.00199999999
